$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The footnote/legend block below the holiday table (rows 21-31) is no
# longer needed, so clear its contents while keeping the existing
# formatting (merged cells, styles, row heights, etc.) intact.
$ws.Range("A21:G31").ClearContents() | Out-Null

# Reflect the selection the user had active right after clearing the block.
$ws.Range("A19:G33").Select() | Out-Null
